$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("courses")

# Update the department value in C2 from the old faculty name to "Community Services"
$ws.Range("C2").Value = "Community Services"
